$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New monthly data rows to append (Aug 2025 - Dec 2025)
$data = @(
    @{ row = 81; anio = 2025; mes = 8;  ipc_inquilinos = 129.1380219516939; fecha = 45870; ipc_oficial = 123.6655045012229 },
    @{ row = 82; anio = 2025; mes = 9;  ipc_inquilinos = 128.8470143760332; fecha = 45901; ipc_oficial = 123.3126918873914 },
    @{ row = 83; anio = 2025; mes = 10; ipc_inquilinos = 129.5238590804855; fecha = 45931; ipc_oficial = 124.1619399490035 },
    @{ row = 84; anio = 2025; mes = 11; ipc_inquilinos = 129.9393556563299; fecha = 45962; ipc_oficial = 124.4023520840922 },
    @{ row = 85; anio = 2025; mes = 12; ipc_inquilinos = 130.5049439590907; fecha = 45992; ipc_oficial = 124.8290576052454 }
)

# The D column (fecha) is formatted as a date (yyyy-mm-dd); match the
# number format already applied to the existing date cells (e.g. D80).
$dateFormat = $ws.Range("D80").NumberFormat

foreach ($item in $data) {
    $r = $item.row
    $ws.Cells.Item($r, 1).Value = $item.anio
    $ws.Cells.Item($r, 2).Value = $item.mes
    $ws.Cells.Item($r, 3).Value = $item.ipc_inquilinos
    $ws.Cells.Item($r, 4).Value = $item.fecha
    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat
    $ws.Cells.Item($r, 5).Value = $item.ipc_oficial
}
